# Auto-generated edit script: update LeveProfit market-data columns (H-N)
# across multiple sheets, per upstream scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 328.30768
$ws.Range("I80").Value = 208.16667
$ws.Range("J80").Value = 431.2857
$ws.Range("K80").Value = 624.50001
$ws.Range("L80").Value = 1293.8571
$ws.Range("M80").Value = 373.49999
$ws.Range("N80").Value = -3289.8571
$ws.Range("H83").Value = 328.30768
$ws.Range("I83").Value = 208.16667
$ws.Range("J83").Value = 431.2857
$ws.Range("K83").Value = 1873.50003
$ws.Range("L83").Value = 3881.5713
$ws.Range("M83").Value = 3118.49997
$ws.Range("N83").Value = -13865.5713
$ws.Range("H125").Value = 1510.3334
$ws.Range("I125").Value = 1312.4
$ws.Range("J125").Value = 2500
$ws.Range("K125").Value = 11811.6
$ws.Range("L125").Value = 22500
$ws.Range("M125").Value = -9351.6
$ws.Range("N125").Value = -27420
$ws.Range("H135").Value = 3000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 3000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 27000
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -32070
$ws.Range("H137").Value = 2210.4666
$ws.Range("I137").Value = 823.55554
$ws.Range("J137").Value = 4290.8335
$ws.Range("K137").Value = 2470.66662
$ws.Range("L137").Value = 12872.5005
$ws.Range("M137").Value = 79.33338000000003
$ws.Range("N137").Value = -17972.5005

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2765.5
$ws.Range("I45").Value = 2137.3
$ws.Range("J45").Value = 4336
$ws.Range("K45").Value = 2137.3
$ws.Range("L45").Value = 4336
$ws.Range("M45").Value = -1760.3
$ws.Range("N45").Value = -5090
$ws.Range("H61").Value = 1915.0952
$ws.Range("I61").Value = 1314.2667
$ws.Range("J61").Value = 3417.1667
$ws.Range("K61").Value = 1314.2667
$ws.Range("L61").Value = 3417.1667
$ws.Range("M61").Value = -1102.2667
$ws.Range("N61").Value = -3841.1667
$ws.Range("H74").Value = 5747.1
$ws.Range("I74").Value = 5251
$ws.Range("J74").Value = 6904.6665
$ws.Range("K74").Value = 5251
$ws.Range("L74").Value = 6904.6665
$ws.Range("M74").Value = -4377
$ws.Range("N74").Value = -8652.6665
$ws.Range("H77").Value = 5747.1
$ws.Range("I77").Value = 5251
$ws.Range("J77").Value = 6904.6665
$ws.Range("K77").Value = 26255
$ws.Range("L77").Value = 34523.3325
$ws.Range("M77").Value = -21887
$ws.Range("N77").Value = -43259.3325
$ws.Range("H122").Value = 2531.25
$ws.Range("I122").Value = 2204
$ws.Range("J122").Value = 2727.6
$ws.Range("K122").Value = 6612
$ws.Range("L122").Value = 8182.799999999999
$ws.Range("M122").Value = -4162
$ws.Range("N122").Value = -13082.8
$ws.Range("H136").Value = 1915.0952
$ws.Range("I136").Value = 1314.2667
$ws.Range("J136").Value = 3417.1667
$ws.Range("K136").Value = 3942.800099999999
$ws.Range("L136").Value = 10251.5001
$ws.Range("M136").Value = -1392.800099999999
$ws.Range("N136").Value = -15351.5001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 404
$ws.Range("I16").Value = 508
$ws.Range("J16").Value = 300
$ws.Range("K16").Value = 508
$ws.Range("L16").Value = 300
$ws.Range("M16").Value = -338
$ws.Range("N16").Value = -640
$ws.Range("H18").Value = 500
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 500
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 500
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -1558
$ws.Range("H22").Value = 388.75
$ws.Range("I22").Value = 351.66666
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 351.66666
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -178.66666
$ws.Range("N22").Value = -846
$ws.Range("H94").Value = 885.7143
$ws.Range("I94").Value = 885.7143
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 885.7143
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -434.7143
$ws.Range("N94").ClearContents()
$ws.Range("H99").Value = 333334660
$ws.Range("I99").Value = 500001000
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 500001000
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -499999502
$ws.Range("N99").Value = -4996

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 623.1667
$ws.Range("I2").Value = 875
$ws.Range("J2").Value = 119.5
$ws.Range("K2").Value = 875
$ws.Range("L2").Value = 119.5
$ws.Range("M2").Value = -762
$ws.Range("N2").Value = -345.5
$ws.Range("H5").Value = 962.7143
$ws.Range("I5").Value = 333.33334
$ws.Range("J5").Value = 1434.75
$ws.Range("K5").Value = 333.33334
$ws.Range("L5").Value = 1434.75
$ws.Range("M5").Value = -221.33334
$ws.Range("N5").Value = -1658.75
$ws.Range("H10").Value = 927.8571
$ws.Range("I10").Value = 999
$ws.Range("J10").Value = 750
$ws.Range("K10").Value = 999
$ws.Range("L10").Value = 750
$ws.Range("M10").Value = -860
$ws.Range("N10").Value = -1028
$ws.Range("H11").Value = 1238.5714
$ws.Range("I11").Value = 570
$ws.Range("J11").Value = 2130
$ws.Range("K11").Value = 570
$ws.Range("L11").Value = 2130
$ws.Range("M11").Value = -430
$ws.Range("N11").Value = -2410
$ws.Range("H12").Value = 235
$ws.Range("I12").Value = 187.5
$ws.Range("J12").Value = 330
$ws.Range("K12").Value = 187.5
$ws.Range("L12").Value = 330
$ws.Range("M12").Value = -17.5
$ws.Range("N12").Value = -670
$ws.Range("H13").Value = 498.75
$ws.Range("I13").Value = 400
$ws.Range("J13").Value = 597.5
$ws.Range("K13").Value = 400
$ws.Range("L13").Value = 597.5
$ws.Range("M13").Value = -261
$ws.Range("N13").Value = -875.5
$ws.Range("H14").Value = 1907.5
$ws.Range("I14").Value = 2675
$ws.Range("J14").Value = 1523.75
$ws.Range("K14").Value = 2675
$ws.Range("L14").Value = 1523.75
$ws.Range("M14").Value = -2505
$ws.Range("N14").Value = -1863.75
$ws.Range("H19").Value = 180.85715
$ws.Range("I19").Value = 83.2
$ws.Range("J19").Value = 425
$ws.Range("K19").Value = 83.2
$ws.Range("L19").Value = 425
$ws.Range("M19").Value = 86.8
$ws.Range("N19").Value = -765
$ws.Range("H22").Value = 1327.1428
$ws.Range("I22").Value = 615.8889
$ws.Range("J22").Value = 2607.4
$ws.Range("K22").Value = 615.8889
$ws.Range("L22").Value = 2607.4
$ws.Range("M22").Value = -265.8889
$ws.Range("N22").Value = -3307.4
$ws.Range("H24").Value = 180.85715
$ws.Range("I24").Value = 83.2
$ws.Range("J24").Value = 425
$ws.Range("K24").Value = 83.2
$ws.Range("L24").Value = 425
$ws.Range("M24").Value = 86.8
$ws.Range("N24").Value = -765
$ws.Range("H62").Value = 1416.3334
$ws.Range("I62").Value = 1250
$ws.Range("J62").Value = 1749
$ws.Range("K62").Value = 1250
$ws.Range("L62").Value = 1749
$ws.Range("M62").Value = -626
$ws.Range("N62").Value = -2997
$ws.Range("H65").Value = 1416.3334
$ws.Range("I65").Value = 1250
$ws.Range("J65").Value = 1749
$ws.Range("K65").Value = 6250
$ws.Range("L65").Value = 8745
$ws.Range("M65").Value = -3130
$ws.Range("N65").Value = -14985
$ws.Range("H107").Value = 1345.625
$ws.Range("I107").Value = 521
$ws.Range("J107").Value = 2170.25
$ws.Range("K107").Value = 521
$ws.Range("L107").Value = 2170.25
$ws.Range("M107").Value = 1399
$ws.Range("N107").Value = -6010.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4220.4
$ws.Range("I80").Value = 3959.45
$ws.Range("J80").Value = 5264.2
$ws.Range("K80").Value = 11878.35
$ws.Range("L80").Value = 15792.6
$ws.Range("M80").Value = -10942.35
$ws.Range("N80").Value = -17664.6
$ws.Range("H83").Value = 4220.4
$ws.Range("I83").Value = 3959.45
$ws.Range("J83").Value = 5264.2
$ws.Range("K83").Value = 35635.05
$ws.Range("L83").Value = 47377.8
$ws.Range("M83").Value = -30955.05
$ws.Range("N83").Value = -56737.8
$ws.Range("H86").Value = 582.5
$ws.Range("I86").Value = 582.5
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1747.5
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -561.5
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 582.5
$ws.Range("I89").Value = 582.5
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 5242.5
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 685.5
$ws.Range("N89").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3831.8
$ws.Range("I80").Value = 3663.3333
$ws.Range("J80").Value = 4084.5
$ws.Range("K80").Value = 3663.3333
$ws.Range("L80").Value = 4084.5
$ws.Range("M80").Value = -2665.3333
$ws.Range("N80").Value = -6080.5
$ws.Range("H83").Value = 3831.8
$ws.Range("I83").Value = 3663.3333
$ws.Range("J83").Value = 4084.5
$ws.Range("K83").Value = 18316.6665
$ws.Range("L83").Value = 20422.5
$ws.Range("M83").Value = -13324.6665
$ws.Range("N83").Value = -30406.5
$ws.Range("H102").Value = 3099.4707
$ws.Range("I102").Value = 2599.3333
$ws.Range("J102").Value = 4299.8
$ws.Range("K102").Value = 2599.3333
$ws.Range("L102").Value = 4299.8
$ws.Range("M102").Value = -977.3332999999998
$ws.Range("N102").Value = -7543.8
$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2000
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -6340

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6089.6665
$ws.Range("I122").Value = 6089.6665
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 18268.9995
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -15818.9995
$ws.Range("N122").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 10002
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 10002
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 20004
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -22126
$ws.Range("H84").Value = 10002
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 10002
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 100020
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -110628
$ws.Range("H107").Value = 83334960
$ws.Range("I107").Value = 166666930
$ws.Range("J107").Value = 2999
$ws.Range("K107").Value = 500000790
$ws.Range("L107").Value = 8997
$ws.Range("M107").Value = -499998870
$ws.Range("N107").Value = -12837
$ws.Range("H132").Value = 1472
$ws.Range("I132").Value = 1333.091
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 3999.273
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -1469.273
$ws.Range("N132").Value = -14060
$ws.Range("H135").Value = 87000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 87000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 87000
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -97140
